# Nueva funcionalidad de hoja de ruta
# Permite programar viajes
#
# Applies the Backlog updates:
#   - Marks row 101 ("clase para simplificar el servidor de reportes") back to "no comenzado"
#   - Adds three new backlog rows (111-113) for the new "hoja de ruta" feature and related tasks

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Existing task status update: row 101 goes from "en proceso" back to "no comenzado"
$ws.Range("B101").Value = "no comenzado"

# New backlog items appended at the bottom of the table
$ws.Range("A111").Value = "fc 5247 defa, habilitar para facturar cantidades especificas de cada remito"
$ws.Range("B111").Value = "no comenzado"

$ws.Range("A112").Value = "citi ventas, desarrollar funcionalidad según sistema NMA"
$ws.Range("B112").Value = "no comenzado"

$ws.Range("A113").Value = "implementar funcionalidad de hoja de ruta"
$ws.Range("B113").Value = "en proceso"

# Update the view state to match where the user ended up working
$win = $excel.ActiveWindow
$null = $ws.Range("A106").Select()
try { $win.TopLeftCell = $ws.Range("A22") } catch {}
try { $win.ScrollRow = 22 } catch {}
try { $win.ScrollColumn = 1 } catch {}
